$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: Summary
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")

$wsSummary.Range("B2").Value = 0.5
$wsSummary.Range("C2").Value = 0.5
$wsSummary.Range("D2").Value = 1
$wsSummary.Range("E2").Value = 0.6666666666666666
$wsSummary.Range("F2").Value = 0.8333333333333334
$wsSummary.Range("G2").Value = 0.9629629629629629
$wsSummary.Range("H2").Value = 0.7854455105275708
$wsSummary.Range("I2").Value = 534
$wsSummary.Range("J2").Value = 534
$wsSummary.Range("K2").Value = 0
$wsSummary.Range("L2").Value = 0

# ---------------------------------------------------------------------------
# Sheet 2: Classification Report
# ---------------------------------------------------------------------------
$wsReport = $wb.Worksheets.Item("Classification Report")

$wsReport.Range("B2").Value = 0
$wsReport.Range("C2").Value = 0
$wsReport.Range("D2").Value = 0

$wsReport.Range("B3").Value = 0.5
$wsReport.Range("C3").Value = 1
$wsReport.Range("D3").Value = 0.6666666666666666

$wsReport.Range("B4").Value = 0.5
$wsReport.Range("C4").Value = 0.5
$wsReport.Range("D4").Value = 0.5
$wsReport.Range("E4").Value = 0.5

$wsReport.Range("B5").Value = 0.25
$wsReport.Range("C5").Value = 0.5
$wsReport.Range("D5").Value = 0.3333333333333333

$wsReport.Range("B6").Value = 0.25
$wsReport.Range("C6").Value = 0.5
$wsReport.Range("D6").Value = 0.3333333333333333

# ---------------------------------------------------------------------------
# Sheet 3: Confusion Matrix
# ---------------------------------------------------------------------------
$wsMatrix = $wb.Worksheets.Item("Confusion Matrix")

$wsMatrix.Range("B2").Value = 0
$wsMatrix.Range("C2").Value = 534

$wsMatrix.Range("B3").Value = 0
$wsMatrix.Range("C3").Value = 534
